# Applies the edits described by the diff:
#  - Community service hours: 50 -> 100
#  - Community service days: 60 -> 90
#  - Completion date: May 18, 2022 -> June 17, 2022
#  - License type: driving -> hunting
#  - Suspension start date: March 09, 2022 -> March 19, 2022
#  - Remove the "Administrative License Suspension is terminated..." sentence

$d = $word.ActiveDocument

# 1) Community service hours: "50" -> "100"
$d.Content.Find.Execute("50", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "100", 2)

# 2) Community service days: "60" -> "90"
$d.Content.Find.Execute("60", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "90", 2)

# 3) Completion date: "May 18, 2022" -> "June 17, 2022"
$d.Content.Find.Execute("May 18, 2022", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "June 17, 2022", 2)

# 4) License type: "driving" -> "hunting"
$d.Content.Find.Execute("driving", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "hunting", 2)

# 5) Suspension start date: "March 09, 2022" -> "March 19, 2022"
$d.Content.Find.Execute("license is suspended from March 09, 2022", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "license is suspended from March 19, 2022", 2)

# 6) Remove the trailing "Administrative License Suspension..." sentence entirely.
$d.Content.Find.Execute("The Administrative License Suspension is terminated and the OBMV form 2261 shall issue. ", `
                         $true, $false, $false, $false, $false, `
                         $true, 1, $false, "", 2)
